$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add two new "Dasar" rows (3. and 4.) to the first table, right after the
#    existing row 2 ("Peraturan Pemerintah Republik Indonesia Nomor 12 ...").
# ---------------------------------------------------------------------------
$dasarTable = $d.Tables.Item(1)

$rFonts = '<w:rFonts w:ascii="Liberation Sans" w:eastAsia="Liberation Sans" w:hAnsi="Liberation Sans" w:cs="Liberation Sans"/>'
$rFontsNoEa = '<w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans" w:cs="Liberation Sans"/>'

function New-DasarRow($table, [string]$num, [string]$bodyText) {
    $row = $table.Rows.Add()

    # --- Column 1: empty label cell (carries the stray bookmarkEnd that the
    #     template keeps duplicating on every "Dasar" row). ---------------
    $cell1 = $table.Cell($row.Index, 1)
    $xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkEnd w:id="0"/><w:pPr><w:tabs><w:tab w:val="center" w:pos="1133"/></w:tabs></w:pPr><w:r><w:rPr>$rFonts<w:color w:val="00000A"/><w:sz w:val="24"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $cell1.Range.InsertXML($xml1)

    # --- Column 2: the "N." numbering cell, as two separate runs. -------
    $cell2 = $table.Cell($row.Index, 2)
    $xml2 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr>$rFontsNoEa<w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr>$rFontsNoEa<w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>$num</w:t></w:r><w:r><w:rPr>$rFontsNoEa<w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $cell2.Range.InsertXML($xml2)

    # --- Column 3: the body text, plus a trailing space run. -------------
    $cell3 = $table.Cell($row.Index, 3)
    $xml3 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr>$rFonts<w:color w:val="00000A"/><w:sz w:val="24"/></w:rPr><w:t>$bodyText</w:t></w:r><w:r><w:rPr>$rFonts<w:color w:val="00000A"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $cell3.Range.InsertXML($xml3)

    # Row height matches the other "Dasar" rows (316 twips == 15.8 points).
    $row.Height = 15.8
}

New-DasarRow $dasarTable "3" "Peraturan Menteri Dalam Negeri Nomor 23 Tahun 2020 tentang`nPerencanaan Pembinaan dan Pengawasan Pemerintahan Daerah`nTahun 2021;"
New-DasarRow $dasarTable "4" "Program Kerja Pengawasan Tahunan (PKPT) Inspektorat Daerah`nKabupaten Sidoarjo Tahun 2021;"

# ---------------------------------------------------------------------------
# 2) Drop the trailing space after "ANDJAR SURJADIANTO, S.Sos." in the
#    "Kepada :" table (first row, name column).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("ANDJAR SURJADIANTO, S.Sos. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ANDJAR SURJADIANTO, S.Sos.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the "3. MOCH. ARTFIANDO, SH / Pengendali Teknis" row from the
#    "Kepada :" table.
# ---------------------------------------------------------------------------
$kepadaTable = $d.Tables.Item(2)
for ($i = 1; $i -le $kepadaTable.Rows.Count; $i++) {
    $r = $kepadaTable.Rows.Item($i)
    if ($r.Range.Text -like "*MOCH. ARTFIANDO*") {
        $r.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Re-wrap two sentences in the "Untuk :" table onto a second line.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("melaksanakan Monitoring Pengelolaan Keuangan Desa di desa wilayah Kecamatan Balongbendo Kabupaten Sidoarjo.", `
                         $true, $false, $false, $false, $false, $true, 1, $false, `
                         "melaksanakan Monitoring Pengelolaan Keuangan Desa di desa wilayah`nKecamatan Balongbendo Kabupaten Sidoarjo.", 2) | Out-Null

$d.Content.Find.Execute("Jangka waktu monitoring selama 7 (tujuh) hari kerja pada periode tanggal 1 s.d 9 Februari 2021.", `
                         $true, $false, $false, $false, $false, $true, 1, $false, `
                         "Jangka waktu monitoring selama 7 (tujuh) hari kerja pada periode tanggal 1 s.d 9`nFebruari 2021.", 2) | Out-Null
